# edit.ps1 - applies the two changes described by the commit diff:
#   1. The "Date" paragraph text changes from 2024-02-05 to 2024-02-14.
#   2. Inside the environmental-optima equation (the one that begins
#      "A_ij(t) = h x 1/(sigma*sqrt(2*pi)) exp(...)"), the squared term's
#      numerator changes from
#           ( E_j(t) - H_i )^2
#      to
#           E_j(t) - H_i
#      i.e. the enclosing parenthesis-delimiter and the outer "^2"
#      superscript are removed, leaving the subtraction bare in the
#      fraction numerator (the outer square around the whole fraction
#      is untouched).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Date paragraph: 2024-02-05 -> 2024-02-14
# ------------------------------------------------------------------
$d.Content.Find.Execute("2024-02-05", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-02-14", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Locate the target equation among the document's OMaths.
#    Use a signature built only from BMP characters (math-italic E/H
#    are outside the BMP and awkward to build without Add-Type, so we
#    key on the italic lower-case h, sigma and minus-sign that are
#    unique, together, to this equation) so the lookup does not depend
#    on a brittle fixed collection index.
# ------------------------------------------------------------------
$sigma = [char]0x03C3
$minus = [char]0x2212
$hLower = [char]0x210E

$target = $null
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $candidate = $d.OMaths.Item($i)
    $t = $candidate.Range.Text
    if ($t.Contains($sigma) -and $t.Contains($minus) -and $t.Contains($hLower)) {
        $target = $candidate
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target equation"
}

$mathNs = "http://schemas.openxmlformats.org/officeDocument/2006/math"

# Full replacement for the equation, identical to the original except
# that the numerator of the squared fraction (E_j(t) - H_i) no longer
# has the enclosing parenthesis delimiter or the outer exponent "2".
$newOMathInner = '<m:sSub><m:e><m:acc><m:accPr><m:chr m:val="' + [char]0x0302 + '" /></m:accPr><m:e><m:r><m:t>A</m:t></m:r></m:e></m:acc></m:e><m:sub><m:r><m:t>i</m:t></m:r><m:r><m:t>j</m:t></m:r></m:sub></m:sSub><m:d><m:dPr><m:begChr m:val="(" /><m:endChr m:val=")" /><m:sepChr m:val="" /><m:grow /></m:dPr><m:e><m:r><m:t>t</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:r><m:t>h</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>' + [char]0x00D7 + '</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>' + $sigma + '</m:t></m:r><m:rad><m:radPr><m:degHide m:val="on" /></m:radPr><m:deg /><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>' + [char]0x03C0 + '</m:t></m:r></m:e></m:rad></m:den></m:f><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>exp</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>' + $minus + '</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:r><m:t>1</m:t></m:r></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f><m:sSup><m:e><m:d><m:dPr><m:begChr m:val="(" /><m:endChr m:val=")" /><m:sepChr m:val="" /><m:grow /></m:dPr><m:e><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSub><m:e><m:r><m:t>E</m:t></m:r></m:e><m:sub><m:r><m:t>j</m:t></m:r></m:sub></m:sSub><m:d><m:dPr><m:begChr m:val="(" /><m:endChr m:val=")" /><m:sepChr m:val="" /><m:grow /></m:dPr><m:e><m:r><m:t>t</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>' + $minus + '</m:t></m:r><m:sSub><m:e><m:r><m:t>H</m:t></m:r></m:e><m:sub><m:r><m:t>i</m:t></m:r></m:sub></m:sSub></m:num><m:den><m:r><m:t>' + $sigma + '</m:t></m:r></m:den></m:f></m:e></m:d></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup>'

$newXml = '<m:oMathPara xmlns:m="' + $mathNs + '"><m:oMathParaPr><m:jc m:val="center"/></m:oMathParaPr><m:oMath>' + $newOMathInner + '</m:oMath></m:oMathPara>'

$range = $target.Range.Duplicate
$range.InsertXML($newXml)
